$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column (H), mirroring the style of the existing header row
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
